$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'27.653.63"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +0.08%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'1.843.55"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  -0.08%  "
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'1.002"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  -0.02%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'312.24"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -1.03%  "
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  -0.03%  "
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'0.4267"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  +0.71%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'0.3616"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  -0.47%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.07304"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  +0.54%  "
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'0.8690"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  -2.12%  "
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  +0.33%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'1.861.53"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  +2.06%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'6.540"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  -0.31%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'5.327"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  -0.10%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'0.06969"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  +1.17%  "
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  +0.03%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'79.43"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  +0.67%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'0.000008958"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  +1.19%  "
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  +0.19%  "
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -1.04%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'27.698.74"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  +0.32%  "
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'4.978"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  +0.15%  "
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  -1.90%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'2.085.72"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  +1.68%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'1.978"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +1.26%  "
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'155.17"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  -0.28%  "
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'18.51"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -1.50%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'120.23"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -2.48%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'5.234"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -0.47%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'0.08903"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  -0.27%  "
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'0.7647"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -1.20%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = "'2.960"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  +1.22%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'4.497"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  -1.58%  "
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'1.124"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  +2.96%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'1.002"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  +0.07%  "
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'0.05429"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  +1.13%  "
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  +0.47%  "
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -0.27%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'2.817"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  +0.81%  "
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'0.1662"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  +0.76%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'0.5065"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -0.21%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'6.565"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  -4.18%  "
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'8.401"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  +1.81%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.06544"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  -0.66%  "
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'106.26"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -0.04%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'1.001"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -0.04%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.4626"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -1.75%  "
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'1.631"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +0.05%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'64.38"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -0.05%  "
$c.Style = "Normal"
